$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.768.86'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '2.347.70'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''543.72'
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').Value = '''136.87'
$ws.Range('E6').Value = '  -2.88%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '''0.525'
$ws.Range('E8').Value = '  -4.85%  '
$ws.Range('D9').Value = '2.346.80'
$ws.Range('E9').Value = '  -0.92%  '
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('D11').Value = '''0.158'
$ws.Range('E11').Value = '  +1.94%  '
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('D13').Value = '''0.343'
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('D14').Value = '''24.70'
$ws.Range('E14').Value = '  -2.68%  '
$ws.Range('D15').Value = '2.773.44'
$ws.Range('E15').Value = '  -0.80%  '
$ws.Range('D16').Value = '60.682.62'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('E17').Value = '  -2.04%  '
$ws.Range('D18').Value = '2.354.25'
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('D19').Value = '''10.63'
$ws.Range('E19').Value = '  +0.74%  '
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').Value = '''318.94'
$ws.Range('E21').Value = '  +0.76%  '
$ws.Range('E22').Value = '  -1.99%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').Value = '''63.30'
$ws.Range('E24').Value = '  +0.75%  '
$ws.Range('D25').Value = '''1.67'
$ws.Range('E25').Value = '  -7.39%  '
$ws.Range('D26').Value = '''8.34'
$ws.Range('E26').Value = '  +7.76%  '
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').Value = '''497.12'
$ws.Range('E29').Value = '  -3.72%  '
$ws.Range('E30').Value = '  -3.99%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').Value = '''0.145'
$ws.Range('E31').Value = '  +0.55%  '
$ws.Range('B32').Value = 'PEPE'
$ws.Range('C32').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D32').Value = '0.0₃0857'
$ws.Range('E32').Value = '  -7.14%  '
$ws.Range('E33').Value = '  -1.96%  '
$ws.Range('E34').Value = '  -3.75%  '
$ws.Range('D35').Value = '''1.00'
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('D36').Value = '''4.60'
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('D37').Value = '''0.375'
$ws.Range('E37').Value = '  +0.41%  '
$ws.Range('E38').Value = '  +2.34%  '
$ws.Range('E39').Value = '  +6.55%  '
$ws.Range('D41').Value = '''143.19'
$ws.Range('E41').Value = '  +4.74%  '
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('E43').Value = '  +0.98%  '
$ws.Range('D44').Value = '''143.33'
$ws.Range('E44').Value = '  +3.16%  '
$ws.Range('E45').Value = '  +0.64%  '
$ws.Range('D46').Value = '''2.03'
$ws.Range('E46').Value = '  -8.91%  '
$ws.Range('D47').Value = '''0.0517'
$ws.Range('E47').Value = '  +0.79%  '
$ws.Range('D48').Value = '''19.05'
$ws.Range('E48').Value = '  -6.01%  '
$ws.Range('D49').Value = '''0.567'
$ws.Range('E49').Value = '  -1.34%  '
$ws.Range('D50').Value = '''0.0901'
$ws.Range('E50').Value = '  -1.38%  '
$ws.Range('E51').Value = '  -1.55%  '
